$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $text)
    # Enter the value as a text formula result, then paste-special as
    # values only. This keeps the cell's type as text (shared string)
    # without Excel re-interpreting the numeric-looking text as a number,
    # and without creating/touching any cell styles (NumberFormat, etc).
    $range.Formula = "=""" + $text + """"
    $range.Copy()
    $range.PasteSpecial(-4163)
}

# Order matters for how new entries land in the shared-strings table;
# update right-to-left (E2, D2, C2) to match the expected layout.
Set-TextValue $ws.Range("E2") "2.4"
Set-TextValue $ws.Range("D2") "0.75"
Set-TextValue $ws.Range("C2") "9.5"

$excel.CutCopyMode = $false
$ws.Range("C2").Select()
